$d = $word.ActiveDocument

# Move to the very end of the document
$sel = $word.Selection
$sel.EndKey(6) | Out-Null   # wdStory = 6

# First new paragraph: empty
$sel.TypeParagraph()

# Second new paragraph: the combined git command text
$sel.TypeParagraph()
$sel.LanguageID = 1051  # wdSlovak
$sel.TypeText("git -C backend add . && git -C frontend add . && git add . && git commit -m `"update frontend a backend`" && git push origin main")
